$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 778; this shifts the existing rows 778-880
# down to 779-881 (matching the diff's renumbering of every row in that
# range) and grows the sheet dimension to A1:R881.
$ws.Rows("778").Insert()

# Populate the newly inserted row 778 with the new weekly data point.
$ws.Range("A778").Value = 10
$ws.Range("B778").Value = "Vega Modelo de Temuco"
$ws.Range("C778").Value = "La Araucanía"
$ws.Range("D778").Value = 45142
$ws.Range("E778").Value = 9
$ws.Range("F778").Value = 100112045
$ws.Range("G778").Value = "Zapallo"
$ws.Range("H778").Value = "Camote"
$ws.Range("I778").Value = "1a (guarda)"
$ws.Range("J778").Value = 900
$ws.Range("K778").Value = 500
$ws.Range("L778").Value = 500
$ws.Range("M778").Value = 500
$ws.Range("N778").Value = "$/kilo (volumen en unidades)"
$ws.Range("O778").Value = "Región del Maule"
$ws.Range("P778").Value = 500
$ws.Range("Q778").Value = 1
$ws.Range("R778").Value = "Hortaliza"
